# Nov 10th - Status
# Append a new daily-status entry (10/11/2021) to the bottom of the log.
#   A95 DATE      -> "10/11/2021"
#   B95 DONE      -> "Internal code-walkthrough on testapps"
#   C95 PROGRESS  -> "Updating the notes"                              (re-used text)
#   B96 DONE      -> "OpenMax : APIs ,Tunneling and components mechanism"
#   C96 PROGRESS  -> "Revising on BT,have to start C-DS-OS"
#   B97 DONE      -> "Discussed doubts and progress in the group"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(95, 1).Value = "10/11/2021"
$ws.Cells.Item(95, 2).Value = "Internal code-walkthrough on testapps"
$ws.Cells.Item(95, 3).Value = "Updating the notes"
$ws.Cells.Item(96, 2).Value = "OpenMax : APIs ,Tunneling and components mechanism"
$ws.Cells.Item(96, 3).Value = "Revising on BT,have to start C-DS-OS"
$ws.Cells.Item(97, 2).Value = "Discussed doubts and progress in the group"

# Match the formatting already used throughout column A (date, centered, text
# format) and columns B/C (wrapped/centered) so no new cell-style entries get
# created - same visual style as every other row in the sheet.
$ws.Range("A95").NumberFormat = "@"
$ws.Range("A95").HorizontalAlignment = -4108
$ws.Range("A95").VerticalAlignment = -4108

$ws.Range("B95:C96").HorizontalAlignment = -4108
$ws.Range("B95:C96").VerticalAlignment = -4108

$ws.Range("B97").HorizontalAlignment = -4108
$ws.Range("B97").VerticalAlignment = -4108

# Move the active selection to the last cell touched, like the author did.
$ws.Range("B97").Select() | Out-Null
